# Auto-generated Excel COM-interop script
# Applies cached market-data value updates (currentAveragePrice* / Leve*Profit* columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 560304.4
$ws.Range("J17").Value = 590910.5
$ws.Range("L17").Value = 1772731.5
$ws.Range("N17").Value = -1773067.5
# Row 40
$ws.Range("H40").Value = 2968.25
$ws.Range("J40").Value = 3154.2727
$ws.Range("L40").Value = 3154.2727
$ws.Range("N40").Value = -3504.2727
# Row 86
$ws.Range("H86").Value = 7448178
$ws.Range("I86").Value = 1759.8889
$ws.Range("J86").Value = 11171387
$ws.Range("K86").Value = 1759.8889
$ws.Range("L86").Value = 11171387
$ws.Range("M86").Value = -636.8888999999999
$ws.Range("N86").Value = -11173633
# Row 89
$ws.Range("H89").Value = 7448178
$ws.Range("I89").Value = 1759.8889
$ws.Range("J89").Value = 11171387
$ws.Range("K89").Value = 8799.4445
$ws.Range("L89").Value = 55856935
$ws.Range("M89").Value = -3183.4445
$ws.Range("N89").Value = -55868167
# Row 107
$ws.Range("H107").Value = 396.8
$ws.Range("I107").Value = 396.8
$ws.Range("K107").Value = 396.8
$ws.Range("M107").Value = 1523.2
# Row 116
$ws.Range("H116").Value = 17195544
$ws.Range("I116").Value = 10050011
$ws.Range("K116").Value = 10050011
$ws.Range("M116").Value = -10046569
# Row 127
$ws.Range("H127").Value = 1187.8667
$ws.Range("I127").Value = 324.33334
$ws.Range("J127").Value = 2483.1667
$ws.Range("K127").Value = 973.0000200000001
$ws.Range("L127").Value = 7449.500100000001
$ws.Range("M127").Value = 3986.99998
$ws.Range("N127").Value = -17369.5001
# Row 137
$ws.Range("H137").Value = 4698.3184
$ws.Range("I137").Value = 3772.1667
$ws.Range("J137").Value = 5809.7
$ws.Range("K137").Value = 11316.5001
$ws.Range("L137").Value = 17429.1
$ws.Range("M137").Value = -8766.500100000001
$ws.Range("N137").Value = -22529.1
# Row 138
$ws.Range("H138").Value = 2223.8447
$ws.Range("I138").Value = 1765.875
$ws.Range("J138").Value = 2787.5
$ws.Range("K138").Value = 5297.625
$ws.Range("L138").Value = 8362.5
$ws.Range("M138").Value = -157.625
$ws.Range("N138").Value = -18642.5
# Row 141
$ws.Range("H141").Value = 3670.3809
$ws.Range("I141").Value = 3670.3809
$ws.Range("K141").Value = 11011.1427
$ws.Range("M141").Value = -5831.1427

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10908.12
$ws.Range("I32").Value = 6673.9673
$ws.Range("K32").Value = 6673.9673
$ws.Range("M32").Value = -6386.9673
# Row 45
$ws.Range("H45").Value = 618815.1
$ws.Range("J45").Value = 1833.8
$ws.Range("L45").Value = 1833.8
$ws.Range("N45").Value = -2587.8
# Row 61
$ws.Range("H61").Value = 3803.7856
$ws.Range("I61").Value = 2552.111
$ws.Range("J61").Value = 6056.8
$ws.Range("K61").Value = 2552.111
$ws.Range("L61").Value = 6056.8
$ws.Range("M61").Value = -2340.111
$ws.Range("N61").Value = -6480.8
# Row 74
$ws.Range("H74").Value = 90914504
$ws.Range("I74").Value = 200004450
$ws.Range("K74").Value = 200004450
$ws.Range("M74").Value = -200003576
# Row 77
$ws.Range("H77").Value = 90914504
$ws.Range("I77").Value = 200004450
$ws.Range("K77").Value = 1000022250
$ws.Range("M77").Value = -1000017882
# Row 110
$ws.Range("H110").Value = 52633390
$ws.Range("I110").Value = 66668100
$ws.Range("J110").Value = 3224.5
$ws.Range("K110").Value = 66668100
$ws.Range("L110").Value = 3224.5
$ws.Range("M110").Value = -66666055
$ws.Range("N110").Value = -7314.5
# Row 122
$ws.Range("H122").Value = 4950.64
$ws.Range("I122").Value = 4675.8
$ws.Range("K122").Value = 14027.4
$ws.Range("M122").Value = -11577.4
# Row 136
$ws.Range("H136").Value = 3803.7856
$ws.Range("I136").Value = 2552.111
$ws.Range("J136").Value = 6056.8
$ws.Range("K136").Value = 7656.333
$ws.Range("L136").Value = 18170.4
$ws.Range("M136").Value = -5106.333
$ws.Range("N136").Value = -23270.4
# Row 139
$ws.Range("H139").Value = 289357.5
$ws.Range("J139").Value = 289357.5
$ws.Range("L139").Value = 289357.5
$ws.Range("N139").Value = -299637.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3113.6
$ws.Range("I86").Value = 2436.75
$ws.Range("K86").Value = 2436.75
$ws.Range("M86").Value = -1313.75
# Row 89
$ws.Range("H89").Value = 3113.6
$ws.Range("I89").Value = 2436.75
$ws.Range("K89").Value = 12183.75
$ws.Range("M89").Value = -6567.75
# Row 105
$ws.Range("H105").Value = 2359.6316
$ws.Range("I105").Value = 2133.8462
$ws.Range("J105").Value = 2848.8333
$ws.Range("K105").Value = 2133.8462
$ws.Range("L105").Value = 2848.8333
$ws.Range("M105").Value = -386.8462
$ws.Range("N105").Value = -6342.8333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 2814645
$ws.Range("J6").Value = 5001.5
$ws.Range("L6").Value = 5001.5
$ws.Range("N6").Value = -5227.5
# Row 58
$ws.Range("H58").Value = 2587.8386
$ws.Range("I58").Value = 1625.4546
$ws.Range("J58").Value = 4940.3335
$ws.Range("K58").Value = 1625.4546
$ws.Range("L58").Value = 4940.3335
$ws.Range("M58").Value = -1422.4546
$ws.Range("N58").Value = -5346.3335
# Row 62
$ws.Range("H62").Value = 110657.5
$ws.Range("I62").Value = 4209.1665
$ws.Range("K62").Value = 4209.1665
$ws.Range("M62").Value = -3585.1665
# Row 65
$ws.Range("H65").Value = 110657.5
$ws.Range("I65").Value = 4209.1665
$ws.Range("K65").Value = 21045.8325
$ws.Range("M65").Value = -17925.8325
# Row 132
$ws.Range("H132").Value = 4177.1113
$ws.Range("I132").Value = 4326.846
$ws.Range("J132").Value = 3787.8
$ws.Range("K132").Value = 12980.538
$ws.Range("L132").Value = 11363.4
$ws.Range("M132").Value = -10450.538
$ws.Range("N132").Value = -16423.4
# Row 134
$ws.Range("H134").Value = 5579.6895
$ws.Range("I134").Value = 5021.625
$ws.Range("K134").Value = 15064.875
$ws.Range("M134").Value = -12529.875
# Row 136
$ws.Range("H136").Value = 2587.8386
$ws.Range("I136").Value = 1625.4546
$ws.Range("J136").Value = 4940.3335
$ws.Range("K136").Value = 4876.3638
$ws.Range("L136").Value = 14821.0005
$ws.Range("M136").Value = -2326.3638
$ws.Range("N136").Value = -19921.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 19
$ws.Range("H19").Value = 1465.625
$ws.Range("J19").Value = 1900
$ws.Range("L19").Value = 5700
$ws.Range("N19").Value = -6048
# Row 103
$ws.Range("H103").Value = 862.5
$ws.Range("I103").Value = 325
$ws.Range("J103").Value = 1400
$ws.Range("K103").Value = 975
$ws.Range("L103").Value = 4200
$ws.Range("M103").Value = -96
$ws.Range("N103").Value = -5958
# Row 139
$ws.Range("H139").Value = 1760934.6
$ws.Range("I139").Value = 2090594.2
$ws.Range("J139").Value = 2750
$ws.Range("K139").Value = 6271782.6
$ws.Range("L139").Value = 8250
$ws.Range("M139").Value = -6266642.6
$ws.Range("N139").Value = -18530

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 582.5
$ws.Range("I2").Value = 1220
$ws.Range("K2").Value = 1220
$ws.Range("M2").Value = -1107
# Row 10
$ws.Range("H10").Value = 1678000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1678000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1678000
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -1678338
# Row 102
$ws.Range("H102").Value = 1937.8
$ws.Range("I102").Value = 899.2727
$ws.Range("K102").Value = 899.2727
$ws.Range("M102").Value = 722.7273
# Row 113
$ws.Range("H113").Value = 3132.2
$ws.Range("I113").Value = 2099.4
$ws.Range("J113").Value = 3648.6
$ws.Range("K113").Value = 2099.4
$ws.Range("L113").Value = 3648.6
$ws.Range("M113").Value = 70.59999999999991
$ws.Range("N113").Value = -7988.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 6570.5
$ws.Range("I68").Value = 3760.6667
$ws.Range("K68").Value = 3760.6667
$ws.Range("M68").Value = -3011.6667
# Row 71
$ws.Range("H71").Value = 6570.5
$ws.Range("I71").Value = 3760.6667
$ws.Range("K71").Value = 18803.3335
$ws.Range("M71").Value = -15059.3335
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("M95").ClearContents()
# Row 122
$ws.Range("H122").Value = 6338.8687
$ws.Range("I122").Value = 6349.3335
$ws.Range("K122").Value = 19048.0005
$ws.Range("M122").Value = -16598.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1062.2
$ws.Range("J107").Value = 499
$ws.Range("L107").Value = 1497
$ws.Range("N107").Value = -5337
